$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A408:A419").NumberFormat = "@"

$ws.Range("A408").Value = "2026-02-06"
$ws.Range("B408").Value = "10:13:45"
$ws.Range("C408").Value = "10:00"
$ws.Range("D408").Value = "Bathroom"
$ws.Range("E408").Value = "No Motion"
$ws.Range("F408").Value = "Inactive"
$ws.Range("A409").Value = "2026-02-06"
$ws.Range("B409").Value = "10:13:48"
$ws.Range("C409").Value = "10:00"
$ws.Range("D409").Value = "Bathroom"
$ws.Range("E409").Value = "No Motion"
$ws.Range("F409").Value = "Inactive"
$ws.Range("A410").Value = "2026-02-06"
$ws.Range("B410").Value = "10:13:50"
$ws.Range("C410").Value = "10:00"
$ws.Range("D410").Value = "Bathroom"
$ws.Range("E410").Value = "Motion Detected"
$ws.Range("F410").Value = "Active"
$ws.Range("A411").Value = "2026-02-06"
$ws.Range("B411").Value = "10:14:00"
$ws.Range("C411").Value = "10:00"
$ws.Range("D411").Value = "Bathroom"
$ws.Range("E411").Value = "No Motion"
$ws.Range("F411").Value = "Inactive"
$ws.Range("A412").Value = "2026-02-06"
$ws.Range("B412").Value = "10:14:05"
$ws.Range("C412").Value = "10:00"
$ws.Range("D412").Value = "Bathroom"
$ws.Range("E412").Value = "No Motion"
$ws.Range("F412").Value = "Inactive"
$ws.Range("A413").Value = "2026-02-06"
$ws.Range("B413").Value = "10:14:10"
$ws.Range("C413").Value = "10:00"
$ws.Range("D413").Value = "Bathroom"
$ws.Range("E413").Value = "No Motion"
$ws.Range("F413").Value = "Inactive"
$ws.Range("A414").Value = "2026-02-06"
$ws.Range("B414").Value = "10:14:15"
$ws.Range("C414").Value = "10:00"
$ws.Range("D414").Value = "Bathroom"
$ws.Range("E414").Value = "No Motion"
$ws.Range("F414").Value = "Inactive"
$ws.Range("A415").Value = "2026-02-06"
$ws.Range("B415").Value = "10:14:20"
$ws.Range("C415").Value = "10:00"
$ws.Range("D415").Value = "Bathroom"
$ws.Range("E415").Value = "No Motion"
$ws.Range("F415").Value = "Inactive"
$ws.Range("A416").Value = "2026-02-06"
$ws.Range("B416").Value = "10:14:25"
$ws.Range("C416").Value = "10:00"
$ws.Range("D416").Value = "Bathroom"
$ws.Range("E416").Value = "No Motion"
$ws.Range("F416").Value = "Inactive"
$ws.Range("A417").Value = "2026-02-06"
$ws.Range("B417").Value = "10:14:30"
$ws.Range("C417").Value = "10:00"
$ws.Range("D417").Value = "Bathroom"
$ws.Range("E417").Value = "No Motion"
$ws.Range("F417").Value = "Inactive"
$ws.Range("A418").Value = "2026-02-06"
$ws.Range("B418").Value = "10:14:35"
$ws.Range("C418").Value = "10:00"
$ws.Range("D418").Value = "Bathroom"
$ws.Range("E418").Value = "No Motion"
$ws.Range("F418").Value = "Inactive"
$ws.Range("A419").Value = "2026-02-06"
$ws.Range("B419").Value = "10:14:41"
$ws.Range("C419").Value = "10:00"
$ws.Range("D419").Value = "Bathroom"
$ws.Range("E419").Value = "No Motion"
$ws.Range("F419").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A278:A289").NumberFormat = "@"
$ws.Range("E278:E289").NumberFormat = "@"

$ws.Range("A278").Value = "2026-02-06"
$ws.Range("B278").Value = "10:13:43"
$ws.Range("C278").Value = "10:00"
$ws.Range("D278").Value = "Bathroom"
$ws.Range("E278").Value = "69.5%"
$ws.Range("F278").Value = "Active"
$ws.Range("A279").Value = "2026-02-06"
$ws.Range("B279").Value = "10:13:46"
$ws.Range("C279").Value = "10:00"
$ws.Range("D279").Value = "Bathroom"
$ws.Range("E279").Value = "69.4%"
$ws.Range("F279").Value = "Active"
$ws.Range("A280").Value = "2026-02-06"
$ws.Range("B280").Value = "10:13:53"
$ws.Range("C280").Value = "10:00"
$ws.Range("D280").Value = "Bathroom"
$ws.Range("E280").Value = "69.2%"
$ws.Range("F280").Value = "Active"
$ws.Range("A281").Value = "2026-02-06"
$ws.Range("B281").Value = "10:13:58"
$ws.Range("C281").Value = "10:00"
$ws.Range("D281").Value = "Bathroom"
$ws.Range("E281").Value = "69.0%"
$ws.Range("F281").Value = "Active"
$ws.Range("A282").Value = "2026-02-06"
$ws.Range("B282").Value = "10:14:03"
$ws.Range("C282").Value = "10:00"
$ws.Range("D282").Value = "Bathroom"
$ws.Range("E282").Value = "68.9%"
$ws.Range("F282").Value = "Active"
$ws.Range("A283").Value = "2026-02-06"
$ws.Range("B283").Value = "10:14:09"
$ws.Range("C283").Value = "10:00"
$ws.Range("D283").Value = "Bathroom"
$ws.Range("E283").Value = "68.7%"
$ws.Range("F283").Value = "Active"
$ws.Range("A284").Value = "2026-02-06"
$ws.Range("B284").Value = "10:14:14"
$ws.Range("C284").Value = "10:00"
$ws.Range("D284").Value = "Bathroom"
$ws.Range("E284").Value = "68.5%"
$ws.Range("F284").Value = "Active"
$ws.Range("A285").Value = "2026-02-06"
$ws.Range("B285").Value = "10:14:18"
$ws.Range("C285").Value = "10:00"
$ws.Range("D285").Value = "Bathroom"
$ws.Range("E285").Value = "68.5%"
$ws.Range("F285").Value = "Active"
$ws.Range("A286").Value = "2026-02-06"
$ws.Range("B286").Value = "10:14:24"
$ws.Range("C286").Value = "10:00"
$ws.Range("D286").Value = "Bathroom"
$ws.Range("E286").Value = "68.7%"
$ws.Range("F286").Value = "Active"
$ws.Range("A287").Value = "2026-02-06"
$ws.Range("B287").Value = "10:14:29"
$ws.Range("C287").Value = "10:00"
$ws.Range("D287").Value = "Bathroom"
$ws.Range("E287").Value = "68.7%"
$ws.Range("F287").Value = "Active"
$ws.Range("A288").Value = "2026-02-06"
$ws.Range("B288").Value = "10:14:34"
$ws.Range("C288").Value = "10:00"
$ws.Range("D288").Value = "Bathroom"
$ws.Range("E288").Value = "68.8%"
$ws.Range("F288").Value = "Active"
$ws.Range("A289").Value = "2026-02-06"
$ws.Range("B289").Value = "10:14:39"
$ws.Range("C289").Value = "10:00"
$ws.Range("D289").Value = "Bathroom"
$ws.Range("E289").Value = "68.9%"
$ws.Range("F289").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A278:A289").NumberFormat = "@"

$ws.Range("A278").Value = "2026-02-06"
$ws.Range("B278").Value = "10:13:44"
$ws.Range("C278").Value = "10:00"
$ws.Range("D278").Value = "Bathroom"
$ws.Range("E278").Value = "28.1C"
$ws.Range("F278").Value = "Active"
$ws.Range("A279").Value = "2026-02-06"
$ws.Range("B279").Value = "10:13:47"
$ws.Range("C279").Value = "10:00"
$ws.Range("D279").Value = "Bathroom"
$ws.Range("E279").Value = "28.1C"
$ws.Range("F279").Value = "Active"
$ws.Range("A280").Value = "2026-02-06"
$ws.Range("B280").Value = "10:13:54"
$ws.Range("C280").Value = "10:00"
$ws.Range("D280").Value = "Bathroom"
$ws.Range("E280").Value = "28.1C"
$ws.Range("F280").Value = "Active"
$ws.Range("A281").Value = "2026-02-06"
$ws.Range("B281").Value = "10:13:59"
$ws.Range("C281").Value = "10:00"
$ws.Range("D281").Value = "Bathroom"
$ws.Range("E281").Value = "28.0C"
$ws.Range("F281").Value = "Active"
$ws.Range("A282").Value = "2026-02-06"
$ws.Range("B282").Value = "10:14:04"
$ws.Range("C282").Value = "10:00"
$ws.Range("D282").Value = "Bathroom"
$ws.Range("E282").Value = "28.1C"
$ws.Range("F282").Value = "Active"
$ws.Range("A283").Value = "2026-02-06"
$ws.Range("B283").Value = "10:14:09"
$ws.Range("C283").Value = "10:00"
$ws.Range("D283").Value = "Bathroom"
$ws.Range("E283").Value = "28.1C"
$ws.Range("F283").Value = "Active"
$ws.Range("A284").Value = "2026-02-06"
$ws.Range("B284").Value = "10:14:15"
$ws.Range("C284").Value = "10:00"
$ws.Range("D284").Value = "Bathroom"
$ws.Range("E284").Value = "28.1C"
$ws.Range("F284").Value = "Active"
$ws.Range("A285").Value = "2026-02-06"
$ws.Range("B285").Value = "10:14:19"
$ws.Range("C285").Value = "10:00"
$ws.Range("D285").Value = "Bathroom"
$ws.Range("E285").Value = "28.0C"
$ws.Range("F285").Value = "Active"
$ws.Range("A286").Value = "2026-02-06"
$ws.Range("B286").Value = "10:14:25"
$ws.Range("C286").Value = "10:00"
$ws.Range("D286").Value = "Bathroom"
$ws.Range("E286").Value = "28.1C"
$ws.Range("F286").Value = "Active"
$ws.Range("A287").Value = "2026-02-06"
$ws.Range("B287").Value = "10:14:30"
$ws.Range("C287").Value = "10:00"
$ws.Range("D287").Value = "Bathroom"
$ws.Range("E287").Value = "28.0C"
$ws.Range("F287").Value = "Active"
$ws.Range("A288").Value = "2026-02-06"
$ws.Range("B288").Value = "10:14:35"
$ws.Range("C288").Value = "10:00"
$ws.Range("D288").Value = "Bathroom"
$ws.Range("E288").Value = "28.0C"
$ws.Range("F288").Value = "Active"
$ws.Range("A289").Value = "2026-02-06"
$ws.Range("B289").Value = "10:14:40"
$ws.Range("C289").Value = "10:00"
$ws.Range("D289").Value = "Bathroom"
$ws.Range("E289").Value = "28.1C"
$ws.Range("F289").Value = "Active"

